$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header order for columns B..M plus the data they carry.
# Header row (row 1)
$ws.Range("A1").Value = "Date"
$ws.Range("B1").Value = "posWordPercentage"
$ws.Range("C1").Value = "negWordPercentage"
$ws.Range("D1").Value = "posPhrasePercentage"
$ws.Range("E1").Value = "negPhrasePercentage"
$ws.Range("F1").Value = "ElapsedMs"
$ws.Range("G1").Value = "wordCount"
$ws.Range("H1").Value = "sentenceCount"
$ws.Range("I1").Value = "posWordCount"
$ws.Range("J1").Value = "negWordCount"
$ws.Range("K1").Value = "positivePhraseCount"
$ws.Range("L1").Value = "negativePhraseCount"
$ws.Range("M1").Value = "Method"

# Data row (row 2) - the "total score function" results
$ws.Range("A2").Value = 42605.455138888887
$ws.Range("B2").Value = 68
$ws.Range("C2").Value = 30
$ws.Range("D2").Value = 8
$ws.Range("E2").Value = 91
$ws.Range("F2").Value = 9738
$ws.Range("G2").Value = 4481
$ws.Range("H2").Value = 761
$ws.Range("I2").Value = 152
$ws.Range("J2").Value = 67
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 21
$ws.Range("M2").Value = "Noun"

# Resize columns to fit the new content (matches the bestFit column widths in the diff)
$ws.Columns("A:M").AutoFit() | Out-Null
